$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (A:H) keep their text representation (no scientific
# notation / leading-zero loss) for numeric-looking strings such as the
# account numbers in column C.
$ws.Range("A2:H7").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("B2").Value = "Q251990"
$ws.Range("C2").Value = "007400000313200019604463"
$ws.Range("D2").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E2").Value = "AWB"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 6750
$ws.Range("J2").Value = 675
$ws.Range("K2").Value = 6075

# Row 3
$ws.Range("A3").Value = "NOUBAIL MOHAMMED"
$ws.Range("B3").Value = "IR801997"
$ws.Range("C3").Value = "007400000313200019604463"
$ws.Range("D3").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E3").Value = "AWB"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 6750
$ws.Range("J3").Value = 675
$ws.Range("K3").Value = 6075

# Row 4
$ws.Range("A4").Value = "MOHAMED BADRANE"
$ws.Range("B4").Value = "I83603"
$ws.Range("C4").Value = "225400000805987601012173"
$ws.Range("D4").Value = "KHOURIBGA"
$ws.Range("E4").Value = "CA"
$ws.Range("F4").Value = "Point de vente"
$ws.Range("G4").Value = "605/KHOURIBGA NAHDA"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 7500
$ws.Range("J4").Value = 375
$ws.Range("K4").Value = 7125

# Row 5
$ws.Range("A5").Value = "ZERNAKH ABDELLAH"
$ws.Range("B5").Value = "IB19558"
$ws.Range("C5").Value = "145101211406073828000084"
$ws.Range("D5").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E5").Value = "BP"
$ws.Range("F5").Value = "Point de vente"
$ws.Range("G5").Value = "052/FKIH BEN SALEH/AV1"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 12000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 12000

# Row 6
$ws.Range("A6").Value = "NASIRI HASNAA"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "546576878798989898090090"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "CIH"
$ws.Range("F6").Value = "Logement de fonction"
$ws.Range("G6").Value = "905/LF/TADLA OUARDIGHA ZAYANE"
$ws.Range("H6").Value = "mensuelle"
$ws.Range("I6").Value = 9999.99
$ws.Range("J6").Value = 999.99
$ws.Range("K6").Value = 9000

# Row 7 (new)
$ws.Range("A7").Value = " "
$ws.Range("B7").Value = " "
$ws.Range("C7").Value = " "
$ws.Range("D7").Value = " "
$ws.Range("E7").Value = " "
$ws.Range("F7").Value = " "
$ws.Range("G7").Value = " "
$ws.Range("H7").Value = " "
$ws.Range("I7").Value = 42999.99
$ws.Range("J7").Value = 2724.99
$ws.Range("K7").Value = 40275
